# "menambahkan fitur update profil"
# Reworks the "divisi" table (drops divisi_kode, renames divisi->sort),
# renames user_divisi/absen key columns from nip to id_karyawan/id_divisi,
# and adds a new "karyawan" profile-field list in column D (rows 21-50)
# used by the new "update profil" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- divisi table (column D) loses its divisi_kode row, "divisi" field
#     shifts up and a new "sort" field appears where "divisi" used to be.
$ws.Range("D3").Value = "divisi"
$ws.Range("D4").Value = "sort"

# --- user_divisi (column E) / absen (column F) now key off the employee
#     and division ids instead of the nip.
$ws.Range("E3").Value = "id_karyawan"
$ws.Range("F3").Value = "id_karyawan"
$ws.Range("E4").Value = "id_divisi"

# --- drop the now-empty created_at/updated_at cells that used to live in
#     the status_pegawai table rows (A5:A6, C5:C6).
$ws.Range("A5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("C6").ClearContents()

# --- restyle the header-ish id cells for karyawan (B2) and divisi (D2)
$ws.Range("B2").Interior.Color = 65535      # yellow fill (matches B12:B25)
$ws.Range("D2").Interior.Color = 10498160   # purple fill (matches E4)

# --- B3/D3 lose their highlight now that they're regular fields again
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").Style = "Normal"

# --- highlight the karyawan address fields the same yellow as the rest
#     of the karyawan column (B12:B25)
$ws.Range("B12:B25").Interior.Color = 65535

# --- new column D: list of karyawan profile fields used by the update
#     profile feature (rows 21-50)
$ws.Range("D21").Value = "idkaryawan"
$ws.Range("D22").Value = "nama_lengkap"
$ws.Range("D23").Value = "nama_panggilan"
$ws.Range("D24").Value = "gelar"
$ws.Range("D25").Value = "tem_lahir"
$ws.Range("D26").Value = "tgl_lahir"
$ws.Range("D27").Value = "j_kel"
$ws.Range("D28").Value = "agama"
$ws.Range("D29").Value = "status"
$ws.Range("D30").Value = "username"
$ws.Range("D31").Value = "jalan_no"
$ws.Range("D32").Value = "rt"
$ws.Range("D33").Value = "rw"
$ws.Range("D34").Value = "desa_kel"
$ws.Range("D35").Value = "kecamatan"
$ws.Range("D36").Value = "kota"
$ws.Range("D37").Value = "kode_pos"
$ws.Range("D38").Value = "jalan_no_domisili"
$ws.Range("D39").Value = "rt_domisili"
$ws.Range("D40").Value = "rw_domisili"
$ws.Range("D41").Value = "desa_kel_domisili"
$ws.Range("D42").Value = "kecamatan_domisili"
$ws.Range("D43").Value = "kota_domisili"
$ws.Range("D44").Value = "kode_pos_domisili"
$ws.Range("D45").Value = "email"
$ws.Range("D46").Value = "telepon"
$ws.Range("D47").Value = "ktp"
$ws.Range("D48").Value = "no_kk"
$ws.Range("D49").Value = "foto"
$ws.Range("D50").Value = "fotoLama"

# --- leave the selection where the author left it
$ws.Range("D50").Select()
